$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 174.4359
$ws.Cells.Item(17, 9).Value = 50.333332
$ws.Cells.Item(17, 10).Value = 184.77777
$ws.Cells.Item(17, 11).Value = 150.999996
$ws.Cells.Item(17, 12).Value = 554.33331
$ws.Cells.Item(17, 13).Value = 17.00000399999999
$ws.Cells.Item(17, 14).Value = -890.33331

$ws.Cells.Item(28, 8).Value = 6937.2856
$ws.Cells.Item(28, 9).Value = 87.84614999999999
$ws.Cells.Item(28, 10).Value = 18067.625
$ws.Cells.Item(28, 11).Value = 87.84614999999999
$ws.Cells.Item(28, 12).Value = 18067.625
$ws.Cells.Item(28, 13).Value = 397.15385
$ws.Cells.Item(28, 14).Value = -19037.625

$ws.Cells.Item(32, 8).Value = 859.94116
$ws.Cells.Item(32, 9).Value = 724.75
$ws.Cells.Item(32, 10).Value = 901.53845
$ws.Cells.Item(32, 11).Value = 724.75
$ws.Cells.Item(32, 12).Value = 901.53845
$ws.Cells.Item(32, 13).Value = -398.75
$ws.Cells.Item(32, 14).Value = -1553.53845

$ws.Cells.Item(98, 8).Value = 1172.8462
$ws.Cells.Item(98, 10).Value = 1335.3334
$ws.Cells.Item(98, 12).Value = 1335.3334
$ws.Cells.Item(98, 14).Value = -4331.3334

$ws.Cells.Item(113, 8).Value = 2105.2856
$ws.Cells.Item(113, 9).Value = 1997.4
$ws.Cells.Item(113, 11).Value = 1997.4
$ws.Cells.Item(113, 13).Value = 1256.6

$ws.Cells.Item(116, 8).Value = 6524.1304
$ws.Cells.Item(116, 9).Value = 7587.5
$ws.Cells.Item(116, 10).Value = 2696
$ws.Cells.Item(116, 11).Value = 7587.5
$ws.Cells.Item(116, 12).Value = 2696
$ws.Cells.Item(116, 13).Value = -4145.5
$ws.Cells.Item(116, 14).Value = -9580

$ws.Cells.Item(122, 8).Value = 1172.8462
$ws.Cells.Item(122, 10).Value = 1335.3334
$ws.Cells.Item(122, 12).Value = 4006.0002
$ws.Cells.Item(122, 14).Value = -8906.0002

$ws.Cells.Item(129, 8).Value = 1047.6604
$ws.Cells.Item(129, 9).Value = 657.3
$ws.Cells.Item(129, 10).Value = 1138.4419
$ws.Cells.Item(129, 11).Value = 1971.9
$ws.Cells.Item(129, 12).Value = 3415.3257
$ws.Cells.Item(129, 13).Value = 3028.1
$ws.Cells.Item(129, 14).Value = -13415.3257

$ws.Cells.Item(133, 8).Value = 100580
$ws.Cells.Item(133, 10).Value = 100580
$ws.Cells.Item(133, 12).Value = 100580
$ws.Cells.Item(133, 14).Value = -110700

$ws.Cells.Item(136, 8).Value = 60780
$ws.Cells.Item(136, 10).Value = 60780
$ws.Cells.Item(136, 12).Value = 60780
$ws.Cells.Item(136, 14).Value = -70980

$ws.Cells.Item(137, 8).Value = 1247
$ws.Cells.Item(137, 9).Value = 1002.7143
$ws.Cells.Item(137, 10).Value = 1817
$ws.Cells.Item(137, 11).Value = 3008.1429
$ws.Cells.Item(137, 12).Value = 5451
$ws.Cells.Item(137, 13).Value = -458.1428999999998
$ws.Cells.Item(137, 14).Value = -10551

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1089.3334
$ws.Cells.Item(2, 9).Value = 1190
$ws.Cells.Item(2, 11).Value = 1190
$ws.Cells.Item(2, 13).Value = -1077

$ws.Cells.Item(32, 8).Value = 406683.12
$ws.Cells.Item(32, 9).Value = 454251.12
$ws.Cells.Item(32, 11).Value = 454251.12
$ws.Cells.Item(32, 13).Value = -453964.12

$ws.Cells.Item(102, 8).Value = 2040
$ws.Cells.Item(102, 9).Value = 2040
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 2040
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).Value = -418
$ws.Cells.Item(102, 14).ClearContents()

$ws.Cells.Item(116, 8).Value = 1089.3334
$ws.Cells.Item(116, 9).Value = 1190
$ws.Cells.Item(116, 11).Value = 1190
$ws.Cells.Item(116, 13).Value = 1104

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1089.3334
$ws.Cells.Item(3, 9).Value = 1190
$ws.Cells.Item(3, 11).Value = 1190
$ws.Cells.Item(3, 13).Value = -1076

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1502.5294
$ws.Cells.Item(31, 9).Value = 1107.9546
$ws.Cells.Item(31, 10).Value = 1801.862
$ws.Cells.Item(31, 11).Value = 1107.9546
$ws.Cells.Item(31, 12).Value = 1801.862
$ws.Cells.Item(31, 13).Value = -812.9546
$ws.Cells.Item(31, 14).Value = -2391.862

$ws.Cells.Item(34, 8).Value = 1502.5294
$ws.Cells.Item(34, 9).Value = 1107.9546
$ws.Cells.Item(34, 10).Value = 1801.862
$ws.Cells.Item(34, 11).Value = 1107.9546
$ws.Cells.Item(34, 12).Value = 1801.862
$ws.Cells.Item(34, 13).Value = -905.9546
$ws.Cells.Item(34, 14).Value = -2205.862

$ws.Cells.Item(122, 8).Value = 1618.3513
$ws.Cells.Item(122, 9).Value = 1372.0625
$ws.Cells.Item(122, 10).Value = 1806
$ws.Cells.Item(122, 11).Value = 4116.1875
$ws.Cells.Item(122, 12).Value = 5418
$ws.Cells.Item(122, 13).Value = -1666.1875
$ws.Cells.Item(122, 14).Value = -10318

$ws.Cells.Item(134, 8).Value = 1199.1666
$ws.Cells.Item(134, 9).Value = 939
$ws.Cells.Item(134, 10).Value = 2500
$ws.Cells.Item(134, 11).Value = 2817
$ws.Cells.Item(134, 12).Value = 7500
$ws.Cells.Item(134, 13).Value = -282
$ws.Cells.Item(134, 14).Value = -12570

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(36, 8).Value = 1000
$ws.Cells.Item(36, 9).Value = 0
$ws.Cells.Item(36, 11).Value = 0
$ws.Cells.Item(36, 13).ClearContents()

$ws.Cells.Item(68, 8).Value = 1074.8429
$ws.Cells.Item(68, 9).Value = 1049
$ws.Cells.Item(68, 10).Value = 1079.15
$ws.Cells.Item(68, 11).Value = 3147
$ws.Cells.Item(68, 12).Value = 3237.45
$ws.Cells.Item(68, 13).Value = -2336
$ws.Cells.Item(68, 14).Value = -4859.450000000001

$ws.Cells.Item(71, 8).Value = 1074.8429
$ws.Cells.Item(71, 9).Value = 1049
$ws.Cells.Item(71, 10).Value = 1079.15
$ws.Cells.Item(71, 11).Value = 9441
$ws.Cells.Item(71, 12).Value = 9712.35
$ws.Cells.Item(71, 13).Value = -5385
$ws.Cells.Item(71, 14).Value = -17824.35

$ws.Cells.Item(80, 8).Value = 2725
$ws.Cells.Item(80, 9).Value = 2633.3333
$ws.Cells.Item(80, 10).Value = 3000
$ws.Cells.Item(80, 11).Value = 7899.999899999999
$ws.Cells.Item(80, 12).Value = 9000
$ws.Cells.Item(80, 13).Value = -6963.999899999999
$ws.Cells.Item(80, 14).Value = -10872

$ws.Cells.Item(83, 8).Value = 2725
$ws.Cells.Item(83, 9).Value = 2633.3333
$ws.Cells.Item(83, 10).Value = 3000
$ws.Cells.Item(83, 11).Value = 23699.9997
$ws.Cells.Item(83, 12).Value = 27000
$ws.Cells.Item(83, 13).Value = -19019.9997
$ws.Cells.Item(83, 14).Value = -36360

$ws.Cells.Item(107, 8).Value = 1308.6
$ws.Cells.Item(107, 9).Value = 258.83783
$ws.Cells.Item(107, 10).Value = 2330.7368
$ws.Cells.Item(107, 11).Value = 776.51349
$ws.Cells.Item(107, 12).Value = 6992.2104
$ws.Cells.Item(107, 13).Value = 1143.48651
$ws.Cells.Item(107, 14).Value = -10832.2104

$ws.Cells.Item(129, 8).Value = 1070
$ws.Cells.Item(129, 9).Value = 655
$ws.Cells.Item(129, 10).Value = 1900
$ws.Cells.Item(129, 11).Value = 1965
$ws.Cells.Item(129, 12).Value = 5700
$ws.Cells.Item(129, 13).Value = 3035
$ws.Cells.Item(129, 14).Value = -15700

$ws.Cells.Item(131, 8).Value = 1058.1464
$ws.Cells.Item(131, 9).Value = 834.3077
$ws.Cells.Item(131, 10).Value = 1162.0714
$ws.Cells.Item(131, 11).Value = 2502.9231
$ws.Cells.Item(131, 12).Value = 3486.2142
$ws.Cells.Item(131, 13).Value = 2537.0769
$ws.Cells.Item(131, 14).Value = -13566.2142

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 3000
$ws.Cells.Item(126, 9).Value = 3000
$ws.Cells.Item(126, 11).Value = 9000
$ws.Cells.Item(126, 13).Value = -6530

$ws.Cells.Item(132, 8).Value = 4499.25
$ws.Cells.Item(132, 9).Value = 2999.5
$ws.Cells.Item(132, 11).Value = 8998.5
$ws.Cells.Item(132, 13).Value = -6468.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 29415326
$ws.Cells.Item(40, 9).Value = 52633580
$ws.Cells.Item(40, 11).Value = 52633580
$ws.Cells.Item(40, 13).Value = -52633444

$ws.Cells.Item(61, 8).Value = 4271.706
$ws.Cells.Item(61, 9).Value = 3920.4443
$ws.Cells.Item(61, 11).Value = 3920.4443
$ws.Cells.Item(61, 13).Value = -3718.4443

$ws.Cells.Item(82, 8).Value = 1832.5
$ws.Cells.Item(82, 9).Value = 1765
$ws.Cells.Item(82, 11).Value = 1765
$ws.Cells.Item(82, 13).Value = -1404

$ws.Cells.Item(85, 8).Value = 1832.5
$ws.Cells.Item(85, 9).Value = 1765
$ws.Cells.Item(85, 11).Value = 1765
$ws.Cells.Item(85, 13).Value = -517

$ws.Cells.Item(93, 8).Value = 11218.363
$ws.Cells.Item(93, 9).Value = 17648.334
$ws.Cells.Item(93, 10).Value = 3502.4
$ws.Cells.Item(93, 11).Value = 17648.334
$ws.Cells.Item(93, 12).Value = 3502.4
$ws.Cells.Item(93, 13).Value = -16400.334
$ws.Cells.Item(93, 14).Value = -5998.4

$ws.Cells.Item(113, 8).Value = 4271.706
$ws.Cells.Item(113, 9).Value = 3920.4443
$ws.Cells.Item(113, 11).Value = 3920.4443
$ws.Cells.Item(113, 13).Value = -1750.4443

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 818.2308
$ws.Cells.Item(107, 9).Value = 670.2
$ws.Cells.Item(107, 11).Value = 2010.6
$ws.Cells.Item(107, 13).Value = -90.60000000000014

$ws.Cells.Item(126, 8).Value = 1502.1111
$ws.Cells.Item(126, 9).Value = 1000.8
$ws.Cells.Item(126, 11).Value = 3002.4
$ws.Cells.Item(126, 13).Value = -532.3999999999996
